$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.87
$ws.Range("K2").Value = 1.92
$ws.Range("L2").Value = 4.5
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 1.82
$ws.Range("R2").Value = 1.92
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 1.25
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 1.17
$ws.Range("AA2").Value = 2.05
$ws.Range("AB2").Value = 1.7
$ws.Range("AD2").Value = 9
$ws.Range("AE2").Value = 9.5
$ws.Range("AF2").Value = 19
$ws.Range("AK2").Value = 17
$ws.Range("AN2").Value = 8.5
$ws.Range("AP2").Value = 15

# Row 3
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 4.2
$ws.Range("K3").Value = 1.77
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 5
$ws.Range("S3").Value = 3.4
$ws.Range("T3").Value = 1.33
$ws.Range("U3").Value = 5.8
$ws.Range("V3").Value = 1.14
$ws.Range("Y3").Value = 1.73
$ws.Range("Z3").Value = 2.08
$ws.Range("AE3").Value = 11
$ws.Range("AN3").Value = 7.5

# Row 4
$ws.Range("K4").Value = 1.87
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 1.94
$ws.Range("R4").Value = 1.79
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 1.14

# Row 5
$ws.Range("G5").Value = 2.9
$ws.Range("H5").Value = 2.75
$ws.Range("I5").Value = 2.9
$ws.Range("K5").Value = 1.69
$ws.Range("U5").Value = 7.2
$ws.Range("AA5").Value = 2.75
$ws.Range("AB5").Value = 1.4
$ws.Range("AG5").Value = 41
$ws.Range("AI5").Value = 4.33
$ws.Range("AR5").Value = 41

# Row 6
$ws.Range("O6").Value = 1.73
$ws.Range("P6").Value = 2
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.3
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 1.08
$ws.Range("Y6").Value = 1.75
$ws.Range("Z6").Value = 2.05

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 2.75
$ws.Range("J7").Value = 3.25
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("AH7").Value = 51
$ws.Range("AO7").Value = 15

# Row 8
$ws.Range("G8").Value = 2.38
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 3.1
$ws.Range("S8").Value = 2.08
$ws.Range("T8").Value = 1.73
$ws.Range("AH8").Value = 29
$ws.Range("AI8").Value = 9
$ws.Range("AJ8").Value = 6.5
$ws.Range("AM8").Value = 251
$ws.Range("AO8").Value = 15

# Row 12
$ws.Range("M12").Value = 1.03
$ws.Range("O12").Value = 1.22
$ws.Range("X12").Value = 1.33

# Row 13
$ws.Range("K13").Value = 1.87
$ws.Range("Q13").Value = 1.97
$ws.Range("R13").Value = 1.77

# Row 14
$ws.Range("S14").Value = 2.15
$ws.Range("T14").Value = 1.63
